$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each year occupies a 4-row block (A,B,C,D quarters) starting at row 2.
# The edit swaps the "B" and "C" quarter rows within every block.
# Column C ("...比上年同期增减") is blank for the older blocks (rows 2-53) and
# only carries real numbers from the 2016 block onward (rows 54-69); leaving
# truly-blank C cells untouched keeps them as the explicit empty cells the
# sheet already had instead of wiping them out via an empty-string write.
for ($blockStart = 2; $blockStart -le 69; $blockStart += 4) {
    $rowB = $blockStart + 1
    $rowC = $blockStart + 2

    $rangeB = $ws.Range("A" + $rowB + ":B" + $rowB)
    $rangeC = $ws.Range("A" + $rowC + ":B" + $rowC)
    $valB = $rangeB.Value()
    $valC = $rangeC.Value()
    $rangeB.Value = $valC
    $rangeC.Value = $valB

    $rangeB2 = $ws.Range("D" + $rowB + ":E" + $rowB)
    $rangeC2 = $ws.Range("D" + $rowC + ":E" + $rowC)
    $valB2 = $rangeB2.Value()
    $valC2 = $rangeC2.Value()
    $rangeB2.Value = $valC2
    $rangeC2.Value = $valB2

    $cB = $ws.Range("C" + $rowB)
    $cC = $ws.Range("C" + $rowC)
    if ($cB.Value() -ne "" -or $cC.Value() -ne "") {
        $cvB = $cB.Value()
        $cvC = $cC.Value()
        $cB.Value = $cvC
        $cC.Value = $cvB
    }
}

# Remove the now-redundant F (产销率) and G (销售量) columns entirely.
$ws.Range("F:G").EntireColumn.Delete()
